# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-28 12:37:14
#
# Normalizes the "Recorded By" (column G) cell values on the "Session Analysis
# Results" sheet: each cell holds a comma-separated list of recorder names /
# emails. Any entry that is exactly "System" is moved to the end of its list
# (other entries keep their relative order); if no "System" entry is present,
# the list is simply reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ',\s*' | ForEach-Object { $_.Trim() }

    if ($parts.Count -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $others = @($parts | Where-Object { -not $_.Equals("System") })
        $systems = @($parts | Where-Object { $_.Equals("System") })
        $newParts = @($others) + @($systems)
    } else {
        $newParts = @($parts[($parts.Count - 1)..0])
    }

    $newVal = ($newParts -join ", ")

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
